$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (keeps them as text, matching source data).
$priceTextRows = @(5,14,15,17,18,22,24,25,26,27,32,44,48,50)
foreach ($r in $priceTextRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply updated cell values
$ws.Cells.Item(2, 4).Value = "27.675.83"
$ws.Cells.Item(2, 5).Value = "  -0.82%  "
$ws.Cells.Item(3, 4).Value = "1.584.90"
$ws.Cells.Item(3, 5).Value = "  -3.11%  "
$ws.Cells.Item(5, 4).Value = "206.45"
$ws.Cells.Item(5, 5).Value = "  -2.46%  "
$ws.Cells.Item(6, 5).Value = "  -2.89%  "
$ws.Cells.Item(7, 5).Value = "  +0.29%  "
$ws.Cells.Item(8, 5).Value = "  -4.81%  "
$ws.Cells.Item(9, 5).Value = "  -1.29%  "
$ws.Cells.Item(10, 5).Value = "  -3.13%  "
$ws.Cells.Item(11, 5).Value = "  -1.91%  "
$ws.Cells.Item(12, 4).Value = "1.810.74"
$ws.Cells.Item(12, 5).Value = "  -3.05%  "
$ws.Cells.Item(13, 4).Value = "1.611.66"
$ws.Cells.Item(13, 5).Value = "  -1.55%  "
$ws.Cells.Item(14, 4).Value = "3.85"
$ws.Cells.Item(14, 5).Value = "  -4.06%  "
$ws.Cells.Item(15, 4).Value = "0.530"
$ws.Cells.Item(15, 5).Value = "  -5.80%  "
$ws.Cells.Item(16, 4).Value = "27.651.09"
$ws.Cells.Item(16, 5).Value = "  -0.99%  "
$ws.Cells.Item(17, 4).Value = "63.24"
$ws.Cells.Item(17, 5).Value = "  -3.14%  "
$ws.Cells.Item(18, 4).Value = "220.10"
$ws.Cells.Item(18, 5).Value = "  -3.85%  "
$ws.Cells.Item(19, 5).Value = "  -3.72%  "
$ws.Cells.Item(20, 5).Value = "  -5.31%  "
$ws.Cells.Item(22, 4).Value = "4.14"
$ws.Cells.Item(22, 5).Value = "  -5.19%  "
$ws.Cells.Item(23, 5).Value = "  -6.46%  "
$ws.Cells.Item(24, 4).Value = "1.97"
$ws.Cells.Item(24, 5).Value = "  -5.75%  "
$ws.Cells.Item(25, 4).Value = "153.88"
$ws.Cells.Item(25, 5).Value = "  -1.37%  "
$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).Value = "6.76"
$ws.Cells.Item(26, 5).Value = "  -2.76%  "
$ws.Cells.Item(27, 2).Value = "BinanceUSD"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(27, 4).Value = "1.00"
$ws.Cells.Item(27, 5).Value = "  +0.36%  "
$ws.Cells.Item(28, 5).Value = "  -2.85%  "
$ws.Cells.Item(29, 5).Value = "  -4.05%  "
$ws.Cells.Item(30, 5).Value = "  -2.65%  "
$ws.Cells.Item(31, 5).Value = "  -3.48%  "
$ws.Cells.Item(32, 4).Value = "3.21"
$ws.Cells.Item(32, 5).Value = "  -5.84%  "
$ws.Cells.Item(33, 4).Value = "1.386.44"
$ws.Cells.Item(33, 5).Value = "  -1.00%  "
$ws.Cells.Item(34, 5).Value = "  -5.21%  "
$ws.Cells.Item(36, 5).Value = "  -5.41%  "
$ws.Cells.Item(37, 5).Value = "  -0.99%  "
$ws.Cells.Item(38, 5).Value = "  -3.05%  "
$ws.Cells.Item(39, 5).Value = "  -3.42%  "
$ws.Cells.Item(40, 5).Value = "  -3.72%  "
$ws.Cells.Item(42, 5).Value = "  -2.76%  "
$ws.Cells.Item(43, 5).Value = "  -4.18%  "
$ws.Cells.Item(44, 4).Value = "2.18"
$ws.Cells.Item(44, 5).Value = "  +1.71%  "
$ws.Cells.Item(45, 5).Value = "  -3.72%  "
$ws.Cells.Item(46, 5).Value = "  -4.20%  "
$ws.Cells.Item(47, 4).Value = "1.721.85"
$ws.Cells.Item(47, 5).Value = "  -3.04%  "
$ws.Cells.Item(48, 4).Value = "87.92"
$ws.Cells.Item(48, 5).Value = "  -0.88%  "
$ws.Cells.Item(49, 5).Value = "  -2.28%  "
$ws.Cells.Item(50, 4).Value = "0.0972"
$ws.Cells.Item(50, 5).Value = "  -5.12%  "
$ws.Cells.Item(51, 5).Value = "  -0.98%  "
